# -----------------------------------------------------------------------
# groupinfo.docx edit: "1." / "2." numbering is added in front of each
# member's "Name" line, the member's full name is broken up into
# individually spell-checked word runs, the PRN number runs are merged
# back into one run, and a whole second "copied" member block (Name,
# PRN NO, Email ID, Contact No - for Pratibha Dilip Kanere) is appended
# after the existing Pratiksha Shivaji Masal block, ending with a blank
# trailing paragraph. The _GoBack bookmark moves from the first PRN NO
# line down onto the newly-added Contact No line.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

function Find-ParagraphByPrefix($doc, $prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $candidate = $doc.Paragraphs.Item($i)
        if ($candidate.Range.Text.StartsWith($prefix)) {
            return $candidate
        }
    }
    throw "Could not find a paragraph starting with '$prefix'"
}

# -------------------------------------------------------------------
# 1) "Name : Pratiksha Shivaji Masal." paragraph ->
#    "1.Name : Pratiksha Shivaji Masal." with finer-grained runs and
#    spell-check proofErr markers around each name word.
# -------------------------------------------------------------------
$frag3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>1.</w:t></w:r><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Name</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> : </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Pratiksha</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Shivaji</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Masal</w:t></w:r><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
$nameParagraph = Find-ParagraphByPrefix $d "Name :"
$nameParagraph.Range.InsertXML($frag3)

# -------------------------------------------------------------------
# 2) "PRN NO : 220980720" + bookmark + "083" paragraph -> the two
#    number runs are merged into a single "220980720083" run and the
#    _GoBack bookmark is dropped from here (it reappears later).
# -------------------------------------------------------------------
$frag4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">PRN </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>NO :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>220980720083</w:t></w:r></w:p>
'@
$prnParagraph = Find-ParagraphByPrefix $d "PRN NO :"
$prnParagraph.Range.InsertXML($frag4)

# -------------------------------------------------------------------
# 3) Append a second, "copied" member block after the last existing
#    paragraph (Contact No : 9763915630/9766328342): a "2.Name : ..."
#    line, a new PRN NO line, a new Email ID line, a new Contact No
#    line (now carrying the _GoBack bookmark), and a trailing blank
#    paragraph.
# -------------------------------------------------------------------
$fragNew = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Name</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> : </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Pratibha</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Dilip</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Kanere</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">PRN </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>NO :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>220980720081</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Email </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>ID :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>pratibhakanere1234@gmail.com</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Contact No : </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>9975796835</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>
'@
$lastIndex = $d.Paragraphs.Count
$lastParagraph = $d.Paragraphs.Item($lastIndex)
$lastParagraph.Range.InsertParagraphAfter()
$insertionParagraph = $d.Paragraphs.Item($lastIndex + 1)
$insertionParagraph.Range.InsertXML($fragNew)
